# Apply cell value updates per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: direct assignment is safe.
$textUpdates = @{
  "D2" = "67.259.64"
  "E2" = "  -0.23%  "
  "D3" = "2.618.58"
  "E3" = "  -0.04%  "
  "E4" = "  -0.10%  "
  "E5" = "  -0.13%  "
  "E6" = "  -1.30%  "
  "E7" = "  -0.01%  "
  "E8" = "  +2.51%  "
  "D9" = "2.617.59"
  "E9" = "  +0.03%  "
  "E10" = "  -3.27%  "
  "E11" = "  +0.60%  "
  "E12" = "  -0.86%  "
  "E13" = "  -2.42%  "
  "E14" = "  +0.15%  "
  "D15" = "3.095.90"
  "E15" = "  +0.00%  "
  "E16" = "  -4.37%  "
  "D17" = "67.078.58"
  "E17" = "  -0.26%  "
  "D18" = "2.619.70"
  "E18" = "  +0.11%  "
  "E19" = "  +0.56%  "
  "E20" = "  -2.47%  "
  "E21" = "  -4.16%  "
  "E22" = "  -0.52%  "
  "E23" = "  +0.86%  "
  "B24" = "Litecoin"
  "C24" = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
  "E24" = "  +4.55%  "
  "B25" = "Dai"
  "C25" = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
  "E25" = "  -0.07%  "
  "E26" = "  -1.71%  "
  "B28" = "Binance-PegBSC-USD"
  "C28" = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
  "E28" = "  -0.02%  "
  "B29" = "Bittensor"
  "C29" = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
  "E29" = "  -3.60%  "
  "B30" = "PEPE"
  "C30" = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
  "E30" = "  -2.89%  "
  "E31" = "  -3.87%  "
  "E32" = "  -2.18%  "
  "E33" = "  -1.73%  "
  "E34" = "  +0.15%  "
  "E35" = "  -5.69%  "
  "E36" = "  -2.89%  "
  "E37" = "  -2.45%  "
  "E38" = "  +2.19%  "
  "E39" = "  -3.10%  "
  "E40" = "  -0.99%  "
  "E41" = "  -3.51%  "
  "E42" = "  -1.97%  "
  "E43" = "  -2.02%  "
  "E44" = "  -0.27%  "
  "E45" = "  -0.02%  "
  "E46" = "  -0.88%  "
  "E47" = "  -0.46%  "
  "D48" = "0.0₆0289"
  "E48" = "  -2.55%  "
  "E49" = "  -0.52%  "
  "E50" = "  +4.76%  "
  "E51" = "  -0.95%  "
}
foreach ($addr in $textUpdates.Keys) {
  $ws.Range($addr).Value = $textUpdates[$addr]
}

# Numeric-looking strings (e.g. "595.44") must stay text, matching the
# original inlineStr cells -- force text format before assigning, then
# restore the default "Normal" style so no visible formatting changes.
$numericLookingUpdates = @{
  "D5" = "595.44"
  "D6" = "153.14"
  "D14" = "27.67"
  "D19" = "364.24"
  "D21" = "7.45"
  "D23" = "2.08"
  "D24" = "71.22"
  "D25" = "0.999"
  "D26" = "9.99"
  "D28" = "1.00"
  "D29" = "582.29"
  "D30" = "0.0000102"
  "D32" = "7.81"
  "D33" = "1.82"
  "D37" = "4.87"
  "D38" = "157.51"
  "D40" = "0.367"
  "D44" = "41.15"
  "D46" = "16.35"
  "D47" = "157.36"
  "D49" = "3.72"
  "D50" = "21.87"
  "D51" = "0.623"
}
foreach ($addr in $numericLookingUpdates.Keys) {
  $cell = $ws.Range($addr)
  $cell.NumberFormat = "@"
  $cell.Value = $numericLookingUpdates[$addr]
  $cell.Style = "Normal"
}
